$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 1: Bad Drivers (rows 3-6) ---
$ws.Range("C3").Value = 11508
$ws.Range("D3").Value = 96.2

$ws.Range("A4").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.120.1.9"
$ws.Range("B4").Value = 20
$ws.Range("C4").Value = 3622
$ws.Range("D4").Value = 98.5

$ws.Range("A5").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.90.0.2"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 219
$ws.Range("D5").Value = 98.59999999999999

$ws.Range("B6").Value = 41
$ws.Range("C6").Value = 15349

# --- Section 2: Good Drivers (rows 14-30) ---
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.160.0.4"
$ws.Range("B14").Value = 96526
$ws.Range("E14").ClearContents()

$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.230.0.8"
$ws.Range("B15").Value = 328411
$ws.Range("E15").ClearContents()

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.200.0.6"
$ws.Range("B16").Value = 143808
$ws.Range("E16").ClearContents()

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.190.0.4"
$ws.Range("B17").Value = 287148
$ws.Range("E17").ClearContents()

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.250.10.1"
$ws.Range("B18").Value = 69578
$ws.Range("E18").ClearContents()

$ws.Range("E19").ClearContents()

$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.30.0.11"
$ws.Range("B20").Value = 67111
$ws.Range("D20").Value = 100
$ws.Range("E20").ClearContents()

$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.30.4.1"
$ws.Range("B21").Value = 13016
$ws.Range("E21").ClearContents()

$ws.Range("A22").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 23.70.2.3"
$ws.Range("B22").Value = 18721
$ws.Range("E22").Value = "'2024-07-23"

$ws.Range("A23").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.10.0.7"
$ws.Range("B23").Value = 66577
$ws.Range("E23").Value = "'2024-05-09"

$ws.Range("A24").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B24").Value = 14239
$ws.Range("E24").Value = "'2022-05-23"

$ws.Range("A25").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B25").Value = 265400
$ws.Range("E25").Value = "'2022-05-01"

$ws.Range("E26").Value = "'2021-01-19"

$ws.Range("A27").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.0.1.1"
$ws.Range("B27").Value = 15730
$ws.Range("E27").Value = "'2020-09-28"

$ws.Range("A28").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.40.2.2"
$ws.Range("B28").Value = 88435
$ws.Range("D28").Value = 99.90000000000001
$ws.Range("E28").Value = "'2019-08-31"

$ws.Range("B30").Value = 46270
